# Apply updated dSF (column F) values per the commit "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -7
    5  = -10
    6  = 1
    7  = -1
    8  = -2
    9  = -4
    10 = -3
    11 = -7
    12 = -4
    13 = -2
    15 = -2
    16 = -2
    17 = -2
    18 = -3
    19 = 2
    20 = -1
    21 = 5
    22 = 4
    23 = -4
    24 = 3
    25 = 5
    26 = -3
    28 = -1
    29 = 6
    30 = -3
    31 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
